# ---------------------------------------------------------------------------
# "add steering axle, refactor linear interpolation"
#
# Adds a second data table (tyre slip angle -> lateral force, normalized by
# its max via linear interpolation helper row) below the existing rpm/torque
# table, formats the normalized column, titles the existing chart, and adds
# a second line chart plotting the new "tyre slip angle" series. Also
# repositions the original chart slightly to make room for the new one.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new data table -------------------------------------------------------

$ws.Range("A18").Value = "slip angle (x) "
$ws.Range("B18").Value = "lateral force (N)"
$ws.Range("C18").Value = "normalized lateral force (N)"

$slipAngles  = @(0, 5, 7, 10, 15, 20, 25, 30, 35, 40, 45, 50)
$lateralForce = @(0, 4.7, 5.3, 5.7, 5.6, 5.4, 5.2, 5.1, 5, 4.8, 4.7, 4.7)

for ($i = 0; $i -lt $slipAngles.Length; $i++) {
    $row = 19 + $i
    $ws.Cells.Item($row, 1).Value = $slipAngles[$i]
    $ws.Cells.Item($row, 2).Value = $lateralForce[$i]
}

# "max" helper row used by the normalization (linear interpolation base)
$ws.Range("A32").Value = "max"
$ws.Range("B32").Formula = "=MAX(B19:B30)"

# normalized lateral force column - single fill creates one shared formula
$ws.Range("C19:C30").Formula = "=B19/`$B`$32"
$ws.Range("C19:C30").NumberFormat = "0.000"

# widen column A so the new "slip angle (x) " header is readable
$ws.Columns.Item(1).ColumnWidth = 12.33

# --- title the existing (rpm/torque/power) chart, then move it up --------

$co1 = $ws.ChartObjects().Item(1)
$chart1 = $co1.Chart
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "BRZ torque profile"

$co1.Left = 280.0556640625
$co1.Top = 14
$co1.Width = 327.1875
$co1.Height = 216

# --- new chart: tyre slip angle vs normalized lateral force ---------------

$co2 = $ws.ChartObjects().Add(265.0556640625, 275, 327.1875, 216)
$chart2 = $co2.Chart
$chart2.ChartType = 65 # xlLineMarkers, matches the style used by chart1

$chart2.SeriesCollection().NewSeries()
$s2 = $chart2.SeriesCollection(1)
$s2.Name = "=Sheet1!`$C`$18"
$s2.XValues = $ws.Range("A19:A30")
$s2.Values = $ws.Range("C19:C30")

$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "tyre slip angle"
$chart2.Legend.Position = -4107 # xlLegendPositionBottom

$co2.Left = 265.0556640625
$co2.Top = 275
$co2.Width = 327.1875
$co2.Height = 216

# --- selection / scroll position matches where the new data was entered --

$excel.ActiveWindow.ScrollRow = 5
$ws.Range("C19:C30").Select()
